$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks first so re-adding them below does not create duplicate link entries
$ws.Hyperlinks.Delete()

# Update column widths (B: 51 -> 52, D: 28 -> 30 characters)
$ws.Columns.Item(2).ColumnWidth = 51.1666666666667
$ws.Columns.Item(4).ColumnWidth = 29.1666666666667

# Refresh data rows 2-11 with the newly scraped listing (run at 2025-10-10 01:16:43)
# Row 2: マッチングアプリのAIレコメンド構築
$ws.Cells.Item(2, 1).Value = "2025-10-10 01:16:43"
$ws.Cells.Item(2, 2).Value = "マッチングアプリのAIレコメンド構築"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5410515"
$ws.Cells.Item(2, 7).Value = 338
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◇アプリ"

# Row 3: 【相談希望】在庫管理・出品補助ツールの開発に関す
$ws.Cells.Item(3, 1).Value = "2025-10-10 01:16:43"
$ws.Cells.Item(3, 2).Value = "【相談希望】在庫管理・出品補助ツールの開発に関するZoom面談依頼"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5398112"
$ws.Cells.Item(3, 7).Value = 158
$ws.Cells.Item(3, 8).Value = "◆ツール,開発 ◇管理"

# Row 4: 【バックエンド開発】ポータルサイトの予約情報管理
$ws.Cells.Item(4, 1).Value = "2025-10-10 01:16:43"
$ws.Cells.Item(4, 2).Value = "【バックエンド開発】ポータルサイトの予約情報管理システム構築"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5410302"
$ws.Cells.Item(4, 7).Value = 143
$ws.Cells.Item(4, 8).Value = "◆開発 ◇サイト"

# Row 5: 海外仕入れ元サイト→ツールを動かす為のCSVファ
$ws.Cells.Item(5, 1).Value = "2025-10-10 01:16:43"
$ws.Cells.Item(5, 2).Value = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5251319"
$ws.Cells.Item(5, 7).Value = 135
$ws.Cells.Item(5, 8).Value = "◆ツール,スクレイピング ◇サイト"

# Row 6: 【新規教育プラットフォーム開発】ノーコード・ロー
$ws.Cells.Item(6, 1).Value = "2025-10-10 01:16:43"
$ws.Cells.Item(6, 2).Value = "【新規教育プラットフォーム開発】ノーコード・ローコードで構築できる学習アプリ開発パートナー募集!"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5410616"
$ws.Cells.Item(6, 7).Value = 93
$ws.Cells.Item(6, 8).Value = "◆開発 ◇アプリ"

# Row 7: 【急募】クローン作成アプリ開発
$ws.Cells.Item(7, 1).Value = "2025-10-10 01:16:43"
$ws.Cells.Item(7, 2).Value = "【急募】クローン作成アプリ開発"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5409967"
$ws.Cells.Item(7, 7).Value = 93
$ws.Cells.Item(7, 8).Value = "◆開発 ◇アプリ"

# Row 8: 大手クレジットカード企業向け、Google Cl
$ws.Cells.Item(8, 1).Value = "2025-10-10 01:16:43"
$ws.Cells.Item(8, 2).Value = "大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5410520"
$ws.Cells.Item(8, 7).Value = 75
$ws.Cells.Item(8, 8).Value = "◆開発"

# Row 9: 大手クレジットカード企業向け、Google Cl
$ws.Cells.Item(9, 1).Value = "2025-10-10 01:16:43"
$ws.Cells.Item(9, 2).Value = "大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件_ワーカー"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5410523"
$ws.Cells.Item(9, 7).Value = 75
$ws.Cells.Item(9, 8).Value = "◆開発"

# Row 10: 丸太木取り自動計算ソフト試作(製材ライン向け)
$ws.Cells.Item(10, 1).Value = "2025-10-10 01:16:43"
$ws.Cells.Item(10, 2).Value = "丸太木取り自動計算ソフト試作(製材ライン向け)"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5410017"
$ws.Cells.Item(10, 7).Value = 25
$ws.Cells.Item(10, 8).Value = ""

# Row 11: 【急募】FXトレード履歴を基にしたEA作成依頼
$ws.Cells.Item(11, 1).Value = "2025-10-10 01:16:43"
$ws.Cells.Item(11, 2).Value = "【急募】FXトレード履歴を基にしたEA作成依頼"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5410127"
$ws.Cells.Item(11, 7).Value = 18
$ws.Cells.Item(11, 8).Value = ""

# Re-create the URL hyperlinks (added in row order, so relationship ids line up: rId1..rId10)
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5410515")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5398112")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5410302")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5251319")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5410616")
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://www.lancers.jp/work/detail/5409967")
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), "https://www.lancers.jp/work/detail/5410520")
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), "https://www.lancers.jp/work/detail/5410523")
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), "https://www.lancers.jp/work/detail/5410017")
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), "https://www.lancers.jp/work/detail/5410127")

# Make sure the URL cells keep using the workbook's existing "Hyperlink" cell style
$ws.Cells.Item(2, 6).Style = "Hyperlink"
$ws.Cells.Item(3, 6).Style = "Hyperlink"
$ws.Cells.Item(4, 6).Style = "Hyperlink"
$ws.Cells.Item(5, 6).Style = "Hyperlink"
$ws.Cells.Item(6, 6).Style = "Hyperlink"
$ws.Cells.Item(7, 6).Style = "Hyperlink"
$ws.Cells.Item(8, 6).Style = "Hyperlink"
$ws.Cells.Item(9, 6).Style = "Hyperlink"
$ws.Cells.Item(10, 6).Style = "Hyperlink"
$ws.Cells.Item(11, 6).Style = "Hyperlink"
